# Dashboard Addon feature file
# Applies the cell content changes to the "Login" sheet (sheet1):
#  - Column E becomes the "menu" list (previously held in column F)
#  - Column F becomes a new "icon" list
#  - New column G ("table") holds staff-data table headers
#  - New column H ("page") holds "Manage *" navigation labels

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 - headers (menu column slides from F to E; new icon/table/page columns)
$ws.Range("E1").Value = "menu"
$ws.Range("F1").Value = "icon"
$ws.Range("G1").Value = "table"
$ws.Range("H1").Value = "page"

# "page" column body (Manage Program / Manage User / Manage Batch)
$ws.Range("H4").Value = "Manage Program"
$ws.Range("H2").Value = "Manage User"
$ws.Range("H3").Value = "Manage Batch"

# "icon" column body
$ws.Range("F2").Value = "User"
$ws.Range("F3").Value = "Staff"
$ws.Range("F4").Value = "Batches"
$ws.Range("F5").Value = "Programs"

# "table" column body (staff data grid headers)
$ws.Range("G2").Value = "Staff Data"
$ws.Range("G3").Value = "#"
$ws.Range("G4").Value = "First Name"
$ws.Range("G5").Value = "Last Name"
$ws.Range("G6").Value = "Phone"

# "menu" column body (moved from column F into column E)
$ws.Range("E2").Value = "Home"
$ws.Range("E3").Value = "Program"
$ws.Range("E4").Value = "Batch"
$ws.Range("E5").Value = "Class"
$ws.Range("E6").Value = "Logout"

# Old menu column (F) no longer holds "Logout" in row 6
$ws.Range("F6").ClearContents()

# Update selection / active cell to H6 to match the new extent of data
$ws.Range("H6").Select()
